# Upper Level BCIO relationships sheet update:
#  - Insert a new relationship row (BCIOR:000009 "has abstinence period")
#    as row 11, pushing every following row down by one.
#  - Backfill several previously-blank Domain/Range (and a couple of
#    Definition) cells that the new canonical version populates.
#  - Correct the Domain/Range wording for "causally influenced by" and
#    "difference between".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the new row for "has abstinence period" at row 11.
$ws.Rows(11).Insert()

$ws.Range("A11").Value = "BCIOR:000009"
$ws.Range("B11").Value = "has abstinence period"
$ws.Range("D11").Value = "participates in [RO:0000056]"
$ws.Range("E11").Value = "A relation that links abstinence from a behaviour to a temporal region during which this personal attribute is true"
$ws.Range("F11").Value = "specifically dependent continuant"
$ws.Range("G11").Value = "occurrent"

# 2) Fix up "causally influenced by" (row 3) and "difference between" (row 5).
$ws.Range("G3").Value = "behaviour intention"
$ws.Range("F5").Value = "effect estimate"
$ws.Range("G5").Value = "outcome estimate"

# 3) Backfill Domain/Range cells on the rows that shifted down after the
#    insert (original rows 12-39, now rows 13-40).
$ws.Range("G13").Value = "behavioural attribute"   # has behavioural attribute
$ws.Range("G14").Value = "animal"                  # has behavioural companion
$ws.Range("G17").Value = "person"                  # has behavioural target
$ws.Range("F23").Value = "research study"          # has study investigator
$ws.Range("G23").Value = "research study investigator"
$ws.Range("F24").Value = "research study"          # has study sample
$ws.Range("G24").Value = "research study sample"
$ws.Range("G27").Value = "person"                  # is enacted by
$ws.Range("G30").Value = "temporal interval"        # occupies temporal region
$ws.Range("G37").Value = "human life function"      # serves behavioural function
$ws.Range("F39").Value = "mechanism of action"      # through
